# "published all changes, commented" -- refresh the per-state product data
# table (columns B:K, rows 2:52) with the latest published figures, widen
# column B to fit its (now best-fit) header, and leave the selection where
# the author's cursor ended up (E10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New published values for columns B..K (blender_pallet .. skiprope_box),
# keyed by row number (row 1 is the header, data runs rows 2-52).
$data = @(
    @(2,7,1,5,1,4,207,0,453,0,9365),
    @(3,47,1,33,2,22,2029,29,3654,15,67576),
    @(4,29,1,20,1,14,958,0,2144,33,40265),
    @(5,68,2,47,2,32,3077,0,5593,130,98305),
    @(6,357,7,245,11,166,15980,524,29719,622,516165),
    @(7,54,1,37,2,26,2360,106,4471,53,77476),
    @(8,34,1,23,1,16,1195,0,2688,60,47680),
    @(9,7,1,5,1,3,17,0,317,0,8951),
    @(10,10,1,7,1,5,240,0,646,0,13428),
    @(11,207,4,142,6,97,9327,334,17292,407,299251),
    @(12,101,2,70,3,47,4540,139,8233,190,146097),
    @(13,14,1,10,1,7,303,0,826,0,18657),
    @(14,30,1,21,1,14,1108,0,2216,40,42348),
    @(15,18,1,13,1,9,806,0,1321,0,25774),
    @(16,115,2,79,4,54,5140,291,9547,145,166029),
    @(17,63,2,44,2,30,2771,0,4999,102,90596),
    @(18,27,1,19,1,13,1186,0,2102,27,38764),
    @(19,42,1,29,2,20,1678,0,3303,0,59670),
    @(20,42,1,29,2,20,1844,0,3422,0,60382),
    @(21,65,2,44,2,30,2585,0,5346,109,92677),
    @(22,57,1,39,2,27,2408,136,4675,68,81600),
    @(23,13,1,9,1,6,501,0,1003,0,18467),
    @(24,92,2,64,3,43,4103,40,7411,140,132866),
    @(25,53,1,37,2,25,2209,92,4123,46,75787),
    @(26,57,1,39,2,27,2464,138,4715,69,81838),
    @(27,27,1,19,1,13,1183,0,2100,27,38751),
    @(28,11,1,8,117,5,256,0,623,8392,14937),
    @(29,100,2,69,3,47,4200,120,8025,180,143201),
    @(30,8,1,5,1,4,47,0,578,0,10120),
    @(31,19,1,13,1,9,518,0,1355,0,25979),
    @(32,13,1,9,1,6,523,0,1019,0,18562),
    @(33,86,2,59,3,40,3516,0,6926,103,122907),
    @(34,20,1,14,1,9,656,0,1419,0,28012),
    @(35,30,1,21,1,14,1063,0,2184,39,42156),
    @(36,180,4,123,6,84,7795,31,14943,255,259008),
    @(37,108,2,74,4,51,4828,215,9018,107,155811),
    @(38,38,1,26,2,18,1374,0,2951,0,53408),
    @(39,39,1,27,2,19,1664,0,3124,0,55693),
    @(40,119,3,82,4,56,5233,0,9748,166,171389),
    @(41,11,1,7,1,5,127,0,806,0,14386),
    @(42,50,1,34,2,23,1946,56,4040,28,71140),
    @(43,9,1,6,1,4,183,0,641,0,12145),
    @(44,66,2,45,2,31,2685,0,5383,115,94148),
    @(45,280,5,192,8,130,12332,640,23209,560,404052),
    @(46,32,1,22,1,15,1171,0,2466,50,45098),
    @(47,80,2,55,3,38,3525,0,6593,74,115108),
    @(48,6,1,5,1,3,243,0,238,0,8480),
    @(49,72,2,50,3,34,3064,0,5718,29,103210),
    @(50,55,1,38,2,26,2141,109,4279,55,77974),
    @(51,17,1,12,1,8,463,0,1111,0,23264),
    @(52,6,1,4,1,3,22,0,355,0,7532)
)

foreach ($row in $data) {
    $r = $row[0]
    $arr = New-Object 'object[,]' 1,10
    for ($i = 0; $i -lt 10; $i++) {
        $arr[0,$i] = $row[$i + 1]
    }
    $ws.Range("B$r`:K$r").Value = $arr
}

# Column B now carries longer values -- give it an explicit best-fit width.
$ws.Range("B1").EntireColumn.ColumnWidth = 11.1666667

# Leave the selection on the cell the author last touched.
$ws.Range("E10").Select() | Out-Null
